# Applies the OtherMC.xlsx edit:
#  - inserts a new "id_DK_Central_IndustryH_Biomass" row (value 0) right after
#    "id_DK_Central_BP_Biomass"
#  - updates the "id_DK_Central_BP_Oil" value
#  - inserts a new "id_DK_Central_GT" row (value 2.663131119234357) right
#    after "id_DK_Central_EP"
#  - removes the whole "id_DK_Decentral_*" block (its data is superseded by
#    the aggregated "id_DK_nan_*" rows that now follow directly)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row 7 for id_DK_Central_IndustryH_Biomass (pushes
#    BP_Coal..WS down by one row).
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "id_DK_Central_IndustryH_Biomass"
$ws.Range("B7").Value = 0

# 2) id_DK_Central_BP_Oil (now at row 13 after the insert above) changes
#    value.
$ws.Range("B13").Value = 10.17394480003523

# 3) Insert a new row 18 for id_DK_Central_GT (right after id_DK_Central_EP,
#    now at row 17), pushing HPstandard..WS down by one row.
$ws.Rows.Item(18).Insert()
$ws.Range("A18").Value = "id_DK_Central_GT"
$ws.Range("B18").Value = 2.663131119234357

# 4) Delete the id_DK_Decentral_* block, now occupying rows 24-44
#    (id_DK_Decentral_BH_Biogas .. id_DK_Decentral_SH).
$ws.Range("A24:A44").EntireRow.Delete()

$wb.Save()
